# Swap the locked contestant seat assignment between row 3 and row 4 on the
# "Seat Assignments" sheet: the two rows trade their ID / ContestantID /
# Seat values (a seat swap), and the row 3 "Notes" cell (an empty string)
# travels with the swap, ending up on row 4 instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Seat Assignments")

# Capture original values before overwriting anything. (Use .Value2 here:
# in this host, the bare .Value getter returns a property descriptor
# instead of the cell content; .Value2 returns the real scalar.)
$row3_ID = $ws.Range("A3").Value2
$row3_ContestantID = $ws.Range("C3").Value2
$row3_Notes = $ws.Range("H3").Value2

$row4_ID = $ws.Range("A4").Value2
$row4_ContestantID = $ws.Range("C4").Value2

# Row 3 takes on row 4's ID / ContestantID, and the seat becomes "A3".
$ws.Range("A3").Value = $row4_ID
$ws.Range("C3").Value = $row4_ContestantID
$ws.Range("E3").Value = "A3"

# Row 4 takes on row 3's original ID / ContestantID, and the seat becomes "A2".
$ws.Range("A4").Value = $row3_ID
$ws.Range("C4").Value = $row3_ContestantID
$ws.Range("E4").Value = "A2"

# The (empty) Notes cell moves from row 3 down to row 4.
$ws.Range("H3").ClearContents()
$ws.Range("H4").Value = $row3_Notes
